$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds numeric-looking price strings (e.g. "28.935.25", "326.17") that must
# remain plain text, matching the original inline-string cell type. Assigning such a
# string straight to .Value lets Excel auto-convert it to a number, so we briefly force
# a Text number format, write the value, then restore the cell's original style so no
# visible formatting changes.

$styleD2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.935.25"
$ws.Range("D2").Style = $styleD2
$ws.Range("E2").Value = "  +1.69%  "

$styleD3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.65"
$ws.Range("D3").Style = $styleD3
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("E4").Value = "  -0.44%  "

$styleD5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.17"
$ws.Range("D5").Style = $styleD5
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("E6").Value = "  -0.26%  "

$styleD7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4567"
$ws.Range("D7").Style = $styleD7
$ws.Range("E7").Value = "  +0.26%  "

$styleD8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3895"
$ws.Range("D8").Style = $styleD8
$ws.Range("E8").Value = "  +1.74%  "

$styleD9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07847"
$ws.Range("D9").Style = $styleD9
$ws.Range("E9").Value = "  +0.32%  "

$styleD10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9862"
$ws.Range("D10").Style = $styleD10
$ws.Range("E10").Value = "  +0.11%  "

$styleD11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.85"
$ws.Range("D11").Style = $styleD11
$ws.Range("E11").Value = "  +1.92%  "

$styleD12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.884.12"
$ws.Range("D12").Style = $styleD12
$ws.Range("E12").Value = "  +2.29%  "

$styleD13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.041"
$ws.Range("D13").Style = $styleD13
$ws.Range("E13").Value = "  +2.09%  "

$styleD14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.681"
$ws.Range("D14").Style = $styleD14
$ws.Range("E14").Value = "  +0.84%  "

$styleD15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06949"
$ws.Range("D15").Style = $styleD15
$ws.Range("E15").Value = "  +0.36%  "

$styleD16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.01"
$ws.Range("D16").Style = $styleD16
$ws.Range("E16").Value = "  +1.89%  "

$ws.Range("E17").Value = "  -0.31%  "

$styleD18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009975"
$ws.Range("D18").Style = $styleD18
$ws.Range("E18").Value = "  +0.47%  "

$styleD19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.98"
$ws.Range("D19").Style = $styleD19
$ws.Range("E19").Value = "  +1.87%  "

$styleD20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = $styleD20
$ws.Range("E20").Value = "  -0.20%  "

$styleD21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.934.78"
$ws.Range("D21").Style = $styleD21
$ws.Range("E21").Value = "  +1.71%  "

$styleD22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.285"
$ws.Range("D22").Style = $styleD22
$ws.Range("E22").Value = "  +0.77%  "

$ws.Range("E23").Value = "  +0.72%  "

$styleD24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.138.06"
$ws.Range("D24").Style = $styleD24
$ws.Range("E24").Value = "  +3.32%  "

$styleD25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.055"
$ws.Range("D25").Style = $styleD25
$ws.Range("E25").Value = "  -1.40%  "

$styleD26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.75"
$ws.Range("D26").Style = $styleD26
$ws.Range("E26").Value = "  +1.63%  "

$styleD27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.30"
$ws.Range("D27").Style = $styleD27
$ws.Range("E27").Value = "  +1.18%  "

$styleD28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.899"
$ws.Range("D28").Style = $styleD28
$ws.Range("E28").Value = "  +4.27%  "

$styleD29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.928"
$ws.Range("D29").Style = $styleD29
$ws.Range("E29").Value = "  +2.56%  "

$styleD30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.63"
$ws.Range("D30").Style = $styleD30
$ws.Range("E30").Value = "  +0.24%  "

$styleD31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09307"
$ws.Range("D31").Style = $styleD31
$ws.Range("E31").Value = "  +0.41%  "

$styleD32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9034"
$ws.Range("D32").Style = $styleD32
$ws.Range("E32").Value = "  -0.01%  "

$styleD33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.291"
$ws.Range("D33").Style = $styleD33
$ws.Range("E33").Value = "  +0.37%  "

$styleD34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.322"
$ws.Range("D34").Style = $styleD34
$ws.Range("E34").Value = "  +0.58%  "

$ws.Range("E35").Value = "  -0.43%  "

$styleD36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.201"
$ws.Range("D36").Style = $styleD36
$ws.Range("E36").Value = "  +4.61%  "

$styleD37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05770"
$ws.Range("D37").Style = $styleD37
$ws.Range("E37").Value = "  +2.19%  "

$styleD38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02065"
$ws.Range("D38").Style = $styleD38
$ws.Range("E38").Value = "  +1.54%  "

$styleD39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.003"
$ws.Range("D39").Style = $styleD39
$ws.Range("E39").Value = "  -0.07%  "

$styleD40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.672"
$ws.Range("D40").Style = $styleD40
$ws.Range("E40").Value = "  +1.15%  "

$styleD41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5673"
$ws.Range("D41").Style = $styleD41
$ws.Range("E41").Value = "  +2.15%  "

$styleD42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1773"
$ws.Range("D42").Style = $styleD42
$ws.Range("E42").Value = "  +0.71%  "

$styleD43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.736"
$ws.Range("D43").Style = $styleD43
$ws.Range("E43").Value = "  +1.21%  "

$styleD44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.288"
$ws.Range("D44").Style = $styleD44
$ws.Range("E44").Value = "  +7.82%  "

$styleD45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.88"
$ws.Range("D45").Style = $styleD45
$ws.Range("E45").Value = "  +2.98%  "

$styleD46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5353"
$ws.Range("D46").Style = $styleD46
$ws.Range("E46").Value = "  +2.35%  "

$styleD47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07037"
$ws.Range("D47").Style = $styleD47
$ws.Range("E47").Value = "  -1.32%  "

$styleD48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.843"
$ws.Range("D48").Style = $styleD48
$ws.Range("E48").Value = "  +2.17%  "

$styleD49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.50"
$ws.Range("D49").Style = $styleD49
$ws.Range("E49").Value = "  +0.82%  "

$styleD50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.511"
$ws.Range("D50").Style = $styleD50
$ws.Range("E50").Value = "  +3.58%  "

$styleD51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.069"
$ws.Range("D51").Style = $styleD51
$ws.Range("E51").Value = "  -4.38%  "
